# "Generate Report for handoff"
#
# Context: this localization-status workbook tracks two source files:
#   - 57b8156c-cf9f-4f67-a157-8cfdef53e762.md
#   - 92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.md
#
# A new handoff report was generated for 57b8156c-...md: its status flips
# from "Handed back: in sync with en-US" to "Ready for handoff" and its
# "Latest Handoff Datetime" advances. The two files also swap display
# order (92b3dd88 now listed first, 57b8156c second) across all three
# sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"

$ws.Range("A3").Value = "57b8156c-cf9f-4f67-a157-8cfdef53e762.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.ffb6f841966544fb26d211805f267cd32d2f57d5.zh-cn.xlf"
$ws.Range("D2").Value = "2016-01-11 07:52:34"
$ws.Range("E2").Value = "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.md"
$ws.Range("F2").Value = "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.ffb6f841966544fb26d211805f267cd32d2f57d5.zh-cn.xlf"
$ws.Range("G2").Value = "2016-01-11 07:53:52"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "57b8156c-cf9f-4f67-a157-8cfdef53e762.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "57b8156c-cf9f-4f67-a157-8cfdef53e762.f66315874c1f42410fcb52d24bdccd96ead7e29a.zh-cn.xlf"
$ws.Range("D3").Value = "2016-01-11 07:55:28"
$ws.Range("E3").Value = "57b8156c-cf9f-4f67-a157-8cfdef53e762.md"
$ws.Range("F3").Value = "57b8156c-cf9f-4f67-a157-8cfdef53e762.f66315874c1f42410fcb52d24bdccd96ead7e29a.zh-cn.xlf"
$ws.Range("G3").Value = "2016-01-11 07:53:52"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.ffb6f841966544fb26d211805f267cd32d2f57d5.de-de.xlf"
$ws.Range("D2").Value = "2016-01-11 07:52:53"
$ws.Range("E2").Value = "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.md"
$ws.Range("F2").Value = "92b3dd88-06b7-4bbf-acd7-ac7f81d3b112.ffb6f841966544fb26d211805f267cd32d2f57d5.de-de.xlf"
$ws.Range("G2").Value = "2016-01-11 07:54:26"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "57b8156c-cf9f-4f67-a157-8cfdef53e762.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "57b8156c-cf9f-4f67-a157-8cfdef53e762.f66315874c1f42410fcb52d24bdccd96ead7e29a.de-de.xlf"
$ws.Range("D3").Value = "2016-01-11 07:55:47"
$ws.Range("E3").Value = "57b8156c-cf9f-4f67-a157-8cfdef53e762.md"
$ws.Range("F3").Value = "57b8156c-cf9f-4f67-a157-8cfdef53e762.f66315874c1f42410fcb52d24bdccd96ead7e29a.de-de.xlf"
$ws.Range("G3").Value = "2016-01-11 07:54:26"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"
